$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "COD. SERVICIO" column header to "COD. REDMINE".
# Changing the cell value updates the shared-strings table and, because the
# cell belongs to Table1, the ListObject's column name follows automatically.
$ws.Range("T1").Value = "COD. REDMINE"

# Restore the view: scroll so column P / row 1 is the top-left visible cell
# and select T3 (instead of the previous AM1:XFD1048576 selection).
$excel.ActiveWindow.ScrollColumn = 16
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("T3").Select()
